$d = $word.ActiveDocument

# Merge the split runs in the Title paragraph into a single run
$d.Content.Find.Execute("Sigma Notation: Answers", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Sigma Notation: Answers", 2)

# Merge the split runs in the Author paragraph into a single run
$d.Content.Find.Execute("Ifan Howells-Baines, Mark Toner", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ifan Howells-Baines, Mark Toner", 2)
